$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the entire C:T block gets fully cleared
$fullClearRows = @(3, 13, 19, 25, 31, 37, 42, 52, 58, 64, 70, 76)
foreach ($r in $fullClearRows) {
    $ws.Range("C" + $r + ":T" + $r).ClearContents()
}

# Rows where C:H become a constant (t-test derived) value and I:T get cleared
$partialRows = @(9, 48)
foreach ($r in $partialRows) {
    $ws.Range("C" + $r + ":H" + $r).Value = 0.6296901998555847
    $ws.Range("I" + $r + ":T" + $r).ClearContents()
}
